$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted into the data table at row 322,
# pushing the existing rows 322:341 down to 323:342 (dimension grows to
# A1:R342).
$ws.Rows(322).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(322, 1).Value = 10
$ws.Cells.Item(322, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(322, 3).Value = "La Araucanía"
$ws.Cells.Item(322, 4).Value = 44585
$ws.Cells.Item(322, 5).Value = 9
$ws.Cells.Item(322, 6).Value = 100112032
$ws.Cells.Item(322, 7).Value = "Zapallo italiano"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 125
$ws.Cells.Item(322, 11).Value = 12000
$ws.Cells.Item(322, 12).Value = 12000
$ws.Cells.Item(322, 13).Value = 12000
$ws.Cells.Item(322, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(322, 15).Value = "Región del Maule"
$ws.Cells.Item(322, 16).Value = 200
$ws.Cells.Item(322, 17).Value = 60
$ws.Cells.Item(322, 18).Value = "Hortaliza"
